# Upload data and R script
# - populate the new "ID" column (A) with the GovData360 indicator ids that
#   already back each row's hyperlink
# - add the missing "Link???" cells (+ hyperlinks) for the three WJP rows
#   that didn't have them yet
# - leave the cursor/selection where the author left it (row 19 area, C40)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column A ("ID") values for the rows that now carry an indicator id
# ---------------------------------------------------------------------
$ids = @{
    4  = 290
    5  = 2763
    6  = 291
    7  = 294
    10 = 747
    11 = 476
    12 = 477
    13 = 478
    14 = 468
    15 = 470
    16 = 471
    17 = 472
    18 = 473
    19 = 28828
    22 = 370
    25 = 27881
    26 = 27882
    27 = 27897
    28 = 27914
    29 = 27926
}

foreach ($row in $ids.Keys) {
    $ws.Cells.Item($row, 1).Value = $ids[$row]
}

# A19 sits in the same shaded band as B19 (style index 10) - match the
# shading by lifting the format off that neighbouring shaded cell.
$ws.Range("B19").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A19").Value = 28828

# ---------------------------------------------------------------------
# 2. Fill in the three still-blank "Link???" cells for the World Justice
#    Project rows (C25:C27) and hook up their hyperlinks, matching the
#    existing C28/C29 pattern in the same block.
# ---------------------------------------------------------------------
$links = @{
    25 = "https://govdata360.worldbank.org/indicators/h27881wjp?country=MLI&indicator=27881&viz=line_chart&years=2013,2018"
    26 = "https://govdata360.worldbank.org/indicators/h27882wjp?country=MLI&indicator=27882&viz=line_chart&years=2013,2018"
    27 = "https://govdata360.worldbank.org/indicators/h27897wjp?country=MLI&indicator=27897&viz=line_chart&years=2013,2018"
}

foreach ($row in 25..27) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = "Link???"
    $ws.Hyperlinks.Add($cell, $links[$row])
    # Hyperlinks.Add re-styles the cell with the built-in Hyperlink style;
    # put it back to the shaded "Link???" look used by C28/C29 below it.
    $ws.Range("C28").Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# 3. Restore the selection/scroll position the workbook was saved with.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C40").Select()

Write-Host "Upload data and R script - edits applied"
